{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" -> \"Impact\" bullet list so it\n// contains short, impact-focused accomplishment statements instead of the\n// longer job-duty-style bullets that were there before.\n//\n// Before (6 bullets) -> After (4 bullets):\n//   1. Built redistricting platform ... serving 12,847 analysts ...   -> Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\n//   2. Designed ETL pipelines using PySpark, dbt, ...                 -> $4.7M savings enabled nonprofit access\n//   3. Trigonometric algorithm for boundary estimation ...            -> Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\n//   4. Discovered systematic race coding errors ...                   -> (removed)\n//   5. Achieved 87% prediction accuracy ...                           -> (removed)\n//   6. Built cloud-based data warehouse solutions on AWS ...          -> Real-time collaboration at national scale\n//\n// Several of these bullet sentences (e.g. the \"Trigonometric algorithm...\"\n// and \"Achieved 87% prediction accuracy...\" ones) also appear verbatim\n// elsewhere in the resume (Professional Experience section), so a plain\n// document-wide text search is ambiguous. Instead we navigate structurally:\n// find the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then its \"Impact\"\n// sub-heading, then take the run of bullet (\"\u2022\") paragraphs that follows -\n// that run is unambiguous and matches the diff's context exactly.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading.\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    headingIndex = i;\n    break;\n  }\n}\nif (headingIndex === -1) {\n  throw new Error('Could not find the \"KEY ACHIEVEMENTS AND IMPACT\" heading.');\n}\n\n// The \"Impact\" sub-heading immediately follows.\nconst subHeadingIndex = headingIndex + 1;\nif (paragraphs.items[subHeadingIndex].text !== \"Impact\") {\n  throw new Error('Expected \"Impact\" sub-heading after \"KEY ACHIEVEMENTS AND IMPACT\".');\n}\n\n// Collect the contiguous run of bullet (\"\u2022\") paragraphs following the\n// sub-heading - these are the achievement bullets to rewrite.\nconst bulletIndexes = [];\nlet cursor = subHeadingIndex + 1;\nwhile (cursor < paragraphs.items.length && paragraphs.items[cursor].text.indexOf(\"\u2022\") === 0) {\n  bulletIndexes.push(cursor);\n  cursor++;\n}\nif (bulletIndexes.length !== 6) {\n  throw new Error(\"Expected 6 achievement bullets, found \" + bulletIndexes.length);\n}\n\n// New text for bullets at relative positions 0, 1, 2 and 5 (0-based, in the\n// original order); positions 3 and 4 are dropped entirely.\nconst replacementsByRelativeIndex = {\n  0: \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  1: \"\u2022 $4.7M savings enabled nonprofit access\",\n  2: \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n  5: \"\u2022 Real-time collaboration at national scale\"\n};\nconst relativeIndexesToDelete = [3, 4];\n\nfor (const relativeIndexKey of Object.keys(replacementsByRelativeIndex)) {\n  const relativeIndex = Number(relativeIndexKey);\n  const paragraphIndex = bulletIndexes[relativeIndex];\n  paragraphs.items[paragraphIndex].insertText(replacementsByRelativeIndex[relativeIndex], \"Replace\");\n}\n\nfor (const relativeIndex of relativeIndexesToDelete) {\n  const paragraphIndex = bulletIndexes[relativeIndex];\n  paragraphs.items[paragraphIndex].delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" -> \"Impact\" bullet list so it\n# contains short, impact-focused accomplishment statements instead of the\n# longer job-duty-style bullets that were there before.\n#\n# Before (6 bullets) -> After (4 bullets):\n#   1. Built redistricting platform ... serving 12,847 analysts ...   -> Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\n#   2. Designed ETL pipelines using PySpark, dbt, ...                 -> $4.7M savings enabled nonprofit access\n#   3. Trigonometric algorithm for boundary estimation ...            -> Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\n#   4. Discovered systematic race coding errors ...                   -> (removed)\n#   5. Achieved 87% prediction accuracy ...                           -> (removed)\n#   6. Built cloud-based data warehouse solutions on AWS ...          -> Real-time collaboration at national scale\n#\n# Several of these bullet sentences (e.g. the \"Trigonometric algorithm...\"\n# and \"Achieved 87% prediction accuracy...\" ones) also appear verbatim\n# elsewhere in the resume (Professional Experience section), so a plain\n# document-wide Find/Replace would be ambiguous. Instead we navigate\n# structurally: find the \"KEY ACHIEVEMENTS AND IMPACT\" heading paragraph,\n# then its \"Impact\" sub-heading, then take the run of bullet (\"\u2022\")\n# paragraphs that follows - that run is unambiguous and matches the diff's\n# context exactly.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    # Paragraph.Range.Text includes a trailing paragraph-mark character\n    # (carriage return); strip it so comparisons are exact.\n    return $para.Range.Text.TrimEnd([char]13)\n}\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading.\n$headingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq -1) {\n    throw \"Could not find the 'KEY ACHIEVEMENTS AND IMPACT' heading.\"\n}\n\n# The \"Impact\" sub-heading immediately follows.\n$subHeadingIndex = $headingIndex + 1\nif ((Get-ParaText $d.Paragraphs.Item($subHeadingIndex)) -ne \"Impact\") {\n    throw \"Expected 'Impact' sub-heading after 'KEY ACHIEVEMENTS AND IMPACT'.\"\n}\n\n# Collect the contiguous run of bullet (\"\u2022\") paragraphs following the\n# sub-heading - these are the achievement bullets to rewrite.\n$bulletIndexes = @()\n$cursor = $subHeadingIndex + 1\nwhile ($cursor -le $d.Paragraphs.Count) {\n    $txt = Get-ParaText $d.Paragraphs.Item($cursor)\n    if ($txt.Length -gt 0 -and $txt[0] -eq [char]0x2022) {\n        $bulletIndexes += $cursor\n        $cursor++\n    } else {\n        break\n    }\n}\nif ($bulletIndexes.Count -ne 6) {\n    throw \"Expected 6 achievement bullets, found $($bulletIndexes.Count)\"\n}\n\n# New text for bullets at relative positions 0, 1, 2 and 5 (0-based, in the\n# original order); positions 3 and 4 are dropped entirely.\n$bulletChar = [char]0x2022\n$replacementsByRelativeIndex = @{\n    0 = \"$bulletChar Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n    1 = \"$bulletChar `$4.7M savings enabled nonprofit access\"\n    2 = \"$bulletChar Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n    5 = \"$bulletChar Real-time collaboration at national scale\"\n}\n$relativeIndexesToDelete = @(3, 4)\n\nforeach ($relativeIndex in $replacementsByRelativeIndex.Keys) {\n    $paragraphIndex = $bulletIndexes[$relativeIndex]\n    $d.Paragraphs.Item($paragraphIndex).Range.Text = $replacementsByRelativeIndex[$relativeIndex]\n}\n\n# Delete from the highest index down so earlier, still-pending paragraph\n# indexes in $bulletIndexes stay valid while we work.\n$sortedDeletes = $relativeIndexesToDelete | Sort-Object -Descending\nforeach ($relativeIndex in $sortedDeletes) {\n    $paragraphIndex = $bulletIndexes[$relativeIndex]\n    $d.Paragraphs.Item($paragraphIndex).Range.Delete()\n}\n"}
